$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.408.27"
$ws.Range("E2").Value = "  +5.49%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.050.84"
$ws.Range("E3").Value = "  +3.90%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.16"
$ws.Range("E5").Value = "  +3.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.652"
$ws.Range("E6").Value = "  +2.96%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.47"
$ws.Range("E7").Value = "  +15.58%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +7.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.82"
$ws.Range("E10").Value = "  +2.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0769"
$ws.Range("E11").Value = "  +5.53%  "

# Row 12
$ws.Range("E12").Value = "  +1.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.915"
$ws.Range("E13").Value = "  -2.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.86"
$ws.Range("E14").Value = "  +3.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.37"
$ws.Range("E15").Value = "  +27.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.352.68"
$ws.Range("E16").Value = "  +4.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.58"
$ws.Range("E17").Value = "  +6.62%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.051.95"
$ws.Range("E18").Value = "  +3.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.310.80"
$ws.Range("E19").Value = "  +5.39%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.70"
$ws.Range("E20").Value = "  +3.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0876"
$ws.Range("E21").Value = "  +4.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.46"
$ws.Range("E22").Value = "  +6.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.01"
$ws.Range("E23").Value = "  +3.45%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.64"
$ws.Range("E24").Value = "  +4.57%  "

# Row 25
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  +5.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  +11.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.23"
$ws.Range("E28").Value = "  -1.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.98"
$ws.Range("E29").Value = "  +4.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.116"
$ws.Range("E30").Value = "  +22.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.30"
$ws.Range("E31").Value = "  +9.72%  "

# Row 32
$ws.Range("E32").Value = "  +3.93%  "

# Row 33
$ws.Range("E33").Value = "  +10.72%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  +10.04%  "

# Row 35
$ws.Range("E35").Value = "  +6.59%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +3.72%  "

# Row 37
$ws.Range("E37").Value = "  +5.43%  "

# Row 38
$ws.Range("E38").Value = "  -0.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.05"
$ws.Range("E39").Value = "  +17.70%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  +33.85%  "

# Row 41
$ws.Range("E41").Value = "  +18.12%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  +5.15%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.24"
$ws.Range("E43").Value = "  +2.55%  "

# Row 44
$ws.Range("E44").Value = "  +6.61%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0218"
$ws.Range("E45").Value = "  +4.71%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.09"
$ws.Range("E46").Value = "  +7.88%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.85"
$ws.Range("E47").Value = "  +6.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.00"
$ws.Range("E48").Value = "  +7.56%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.421.42"
$ws.Range("E49").Value = "  +4.23%  "

# Row 50
$ws.Range("E50").Value = "  +2.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.82"
$ws.Range("E51").Value = "  +1.76%  "
